$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This script rewrites the cells in rows 2-21 of the "Artfynd" export sheet
# to match the re-fetched source data (rows were re-ordered/updated upstream).

# Force plain-text format on date-like cells so values are not auto-converted to date serials.
$dateTextCells = @("Y2","AA2","Y3","AA3","Y4","AA4","Y12","AA12","Y13","AA13","Y14","AA14")
foreach ($addr in $dateTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Force text format on new blank placeholder cells so they are materialized.
$blankPlaceholderCells = @("J5","L5","N5","AF5")
foreach ($addr in $blankPlaceholderCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Apply cell value updates ---
$ws.Range("A2").Value = 108415209
$ws.Range("B2").Value = 98520
$ws.Range("E2").Value = 222498
$ws.Range("F2").Value = "Blåsippa"
$ws.Range("G2").Value = "Hepatica nobilis"
$ws.Range("H2").Value = "Schreb."
$ws.Range("P2").Value = "..., Upl"
$ws.Range("Q2").Value = 669035.6441127005
$ws.Range("R2").Value = 6706115.064604912
$ws.Range("Y2").Value = "2023-04-23"
$ws.Range("AA2").Value = "2023-04-23"
$ws.Range("AW2").Value = "Signe Propst"
$ws.Range("AX2").Value = "Signe Propst"
$ws.Range("A3").Value = 108415211
$ws.Range("B3").Value = 98520
$ws.Range("E3").Value = 222498
$ws.Range("F3").Value = "Blåsippa"
$ws.Range("G3").Value = "Hepatica nobilis"
$ws.Range("H3").Value = "Schreb."
$ws.Range("P3").Value = "..., Upl"
$ws.Range("Q3").Value = 668952.4726062249
$ws.Range("R3").Value = 6706084.007164106
$ws.Range("Y3").Value = "2023-04-23"
$ws.Range("AA3").Value = "2023-04-23"
$ws.Range("AW3").Value = "Signe Propst"
$ws.Range("AX3").Value = "Signe Propst"
$ws.Range("A4").Value = 108415210
$ws.Range("B4").Value = 98520
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = "Blåsippa"
$ws.Range("G4").Value = "Hepatica nobilis"
$ws.Range("H4").Value = "Schreb."
$ws.Range("P4").Value = "..., Upl"
$ws.Range("Q4").Value = 668964.6286571589
$ws.Range("R4").Value = 6706131.015868785
$ws.Range("Y4").Value = "2023-04-23"
$ws.Range("AA4").Value = "2023-04-23"
$ws.Range("AW4").Value = "Signe Propst"
$ws.Range("AX4").Value = "Signe Propst"
$ws.Range("A5").Value = 108415170
$ws.Range("B5").Value = 101680
$ws.Range("E5").Value = 222412
$ws.Range("F5").Value = "Tibast"
$ws.Range("G5").Value = "Daphne mezereum"
$ws.Range("H5").Value = "L."
$ws.Range("Q5").Value = 669045.2840127514
$ws.Range("R5").Value = 6706162.943860305
$ws.Range("A6").Value = 108415206
$ws.Range("Q6").Value = 669209.5994803488
$ws.Range("R6").Value = 6706225.957262815
$ws.Range("A7").Value = 108415208
$ws.Range("Q7").Value = 669171.637588384
$ws.Range("R7").Value = 6706286.4308879
$ws.Range("A8").Value = 108415205
$ws.Range("Q8").Value = 669243.9638568971
$ws.Range("R8").Value = 6706189.523024754
$ws.Range("A9").Value = 108415207
$ws.Range("Q9").Value = 669176.7729262316
$ws.Range("R9").Value = 6706271.849970463
$ws.Range("A10").Value = 108415202
$ws.Range("P10").Value = "..., Upl"
$ws.Range("Q10").Value = 669270.4582348154
$ws.Range("R10").Value = 6706056.873020063
$ws.Range("S10").Value = 10
$ws.Range("AW10").Value = "Signe Propst"
$ws.Range("AX10").Value = "Signe Propst"
$ws.Range("A11").Value = 108415204
$ws.Range("Q11").Value = 669262.2281218363
$ws.Range("R11").Value = 6706158.758182475
$ws.Range("A12").Value = 92231854
$ws.Range("B12").Value = 101680
$ws.Range("E12").Value = 222412
$ws.Range("F12").Value = "Tibast"
$ws.Range("G12").Value = "Daphne mezereum"
$ws.Range("H12").Value = "L."
$ws.Range("P12").Value = "Skatenskogen, Upl"
$ws.Range("Q12").Value = 668813.7819832751
$ws.Range("R12").Value = 6705779.112063986
$ws.Range("Y12").Value = "2021-04-02"
$ws.Range("AA12").Value = "2021-04-02"
$ws.Range("AW12").Value = "Ingemar Södergren"
$ws.Range("AX12").Value = "Ingemar Södergren, Barbara Kühn, Julian Klein, Louis Mielke"
$ws.Range("A13").Value = 92231855
$ws.Range("B13").Value = 101680
$ws.Range("E13").Value = 222412
$ws.Range("F13").Value = "Tibast"
$ws.Range("G13").Value = "Daphne mezereum"
$ws.Range("H13").Value = "L."
$ws.Range("P13").Value = "Skatenskogen, Upl"
$ws.Range("Q13").Value = 668764.0189416714
$ws.Range("R13").Value = 6705795.566670171
$ws.Range("Y13").Value = "2021-04-02"
$ws.Range("AA13").Value = "2021-04-02"
$ws.Range("AW13").Value = "Ingemar Södergren"
$ws.Range("AX13").Value = "Ingemar Södergren, Barbara Kühn, Julian Klein, Louis Mielke"
$ws.Range("A14").Value = 92231853
$ws.Range("B14").Value = 101680
$ws.Range("E14").Value = 222412
$ws.Range("F14").Value = "Tibast"
$ws.Range("G14").Value = "Daphne mezereum"
$ws.Range("H14").Value = "L."
$ws.Range("P14").Value = "Skatenskogen, Upl"
$ws.Range("Q14").Value = 668819.8689593319
$ws.Range("R14").Value = 6705775.937335981
$ws.Range("S14").Value = 10
$ws.Range("Y14").Value = "2021-04-02"
$ws.Range("AA14").Value = "2021-04-02"
$ws.Range("AW14").Value = "Ingemar Södergren"
$ws.Range("AX14").Value = "Ingemar Södergren, Barbara Kühn, Julian Klein, Louis Mielke"
$ws.Range("A15").Value = 108409856
$ws.Range("Q15").Value = 668506.5236493029
$ws.Range("R15").Value = 6705660.536876176
$ws.Range("A16").Value = 108409858
$ws.Range("B16").Value = 98520
$ws.Range("E16").Value = 222498
$ws.Range("F16").Value = "Blåsippa"
$ws.Range("G16").Value = "Hepatica nobilis"
$ws.Range("H16").Value = "Schreb."
$ws.Range("Q16").Value = 668388.479731641
$ws.Range("R16").Value = 6705421.322440483
$ws.Range("A17").Value = 108409853
$ws.Range("Q17").Value = 668985.6790587006
$ws.Range("R17").Value = 6705807.392752416
$ws.Range("A18").Value = 108409821
$ws.Range("B18").Value = 95519
$ws.Range("E18").Value = 221945
$ws.Range("F18").Value = "Revlummer"
$ws.Range("G18").Value = "Lycopodium annotinum"
$ws.Range("H18").Value = "L."
$ws.Range("P18").Value = "Skaten, Upl"
$ws.Range("Q18").Value = 668583.1734774174
$ws.Range("R18").Value = 6705693.752984518
$ws.Range("S18").Value = 15
$ws.Range("AW18").Value = "Isac Carlsson"
$ws.Range("AX18").Value = "Isac Carlsson"
$ws.Range("A19").Value = 108409854
$ws.Range("P19").Value = "Skaten, Upl"
$ws.Range("Q19").Value = 668626.4819316415
$ws.Range("R19").Value = 6705752.098212075
$ws.Range("S19").Value = 15
$ws.Range("AW19").Value = "Isac Carlsson"
$ws.Range("AX19").Value = "Isac Carlsson"
$ws.Range("A20").Value = 108415198
$ws.Range("P20").Value = "..., Upl"
$ws.Range("Q20").Value = 668981.6375816888
$ws.Range("R20").Value = 6705904.539491506
$ws.Range("S20").Value = 10
$ws.Range("AW20").Value = "Signe Propst"
$ws.Range("AX20").Value = "Signe Propst"
$ws.Range("A21").Value = 108409852
$ws.Range("B21").Value = 98520
$ws.Range("E21").Value = 222498
$ws.Range("F21").Value = "Blåsippa"
$ws.Range("G21").Value = "Hepatica nobilis"
$ws.Range("H21").Value = "Schreb."
$ws.Range("P21").Value = "Skaten, Upl"
$ws.Range("Q21").Value = 668972.9605770472
$ws.Range("R21").Value = 6705857.196229734
$ws.Range("S21").Value = 15
$ws.Range("AW21").Value = "Isac Carlsson"
$ws.Range("AX21").Value = "Isac Carlsson"
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = "blomning"
$ws.Range("L5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("AF5").Value = ""
$ws.Range("J21").ClearContents()
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("AF21").ClearContents()
